$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared-string rich text cells) ---
# A8: "Volume 32   Number  9" -> "Volume 32   Number  10"
$ws.Range("A8").Value = "Volume 32   Number  10"
# C9: "Report Covering the Week  2/24/2025  Through  3/2/2025"
#  -> "Report Covering the Week  3/3/2025  Through  3/9/2025"
$ws.Range("C9").Value = "Report Covering the Week  3/3/2025  Through  3/9/2025"

# --- Crime Complaints table updates (rows 14-33) ---

# --- Row 14 ---
$ws.Range("F14").Value2 = 1
$ws.Range("N14").Value2 = -82.758620689655

# --- Row 15 ---
$ws.Range("C15").Value2 = 4
$ws.Range("D15").Value2 = 4
$ws.Range("G15").Value2 = 16
$ws.Range("H15").Value2 = -25
$ws.Range("I15").Value2 = 36
$ws.Range("J15").Value2 = 30
$ws.Range("K15").Value2 = 20
$ws.Range("L15").Value2 = -2.702702702702
$ws.Range("M15").Value2 = 71.428571428571
$ws.Range("N15").Value2 = 5.882352941176

# --- Row 16 ---
$ws.Range("C16").Value2 = 21
$ws.Range("D16").Value2 = 35
$ws.Range("E16").Value2 = -40
$ws.Range("F16").Value2 = 93
$ws.Range("G16").Value2 = 171
$ws.Range("H16").Value2 = -45.614035087719
$ws.Range("I16").Value2 = 217
$ws.Range("J16").Value2 = 398
$ws.Range("K16").Value2 = -45.477386934673
$ws.Range("L16").Value2 = -36.549707602339
$ws.Range("M16").Value2 = -41.509433962264
$ws.Range("N16").Value2 = -88.004422332780

# --- Row 17 ---
$ws.Range("C17").Value2 = 50
$ws.Range("D17").Value2 = 46
$ws.Range("E17").Value2 = 8.695652173913
$ws.Range("F17").Value2 = 228
$ws.Range("G17").Value2 = 228
$ws.Range("I17").Value2 = 535
$ws.Range("J17").Value2 = 510
$ws.Range("K17").Value2 = 4.901960784313
$ws.Range("L17").Value2 = 11.924686192468
$ws.Range("M17").Value2 = 122.916666666667
$ws.Range("N17").Value2 = 15.800865800865

# --- Row 18 ---
$ws.Range("C18").Value2 = 33
$ws.Range("D18").Value2 = 38
$ws.Range("E18").Value2 = -13.157894736842
$ws.Range("F18").Value2 = 140
$ws.Range("G18").Value2 = 176
$ws.Range("H18").Value2 = -20.454545454545
$ws.Range("I18").Value2 = 354
$ws.Range("J18").Value2 = 390
$ws.Range("K18").Value2 = -9.230769230769
$ws.Range("L18").Value2 = -15.513126491646
$ws.Range("M18").Value2 = -36.101083032491
$ws.Range("N18").Value2 = -88.188188188188

# --- Row 19 ---
$ws.Range("C19").Value2 = 111
$ws.Range("D19").Value2 = 116
$ws.Range("E19").Value2 = -4.310344827586
$ws.Range("F19").Value2 = 395
$ws.Range("G19").Value2 = 491
$ws.Range("H19").Value2 = -19.551934826883
$ws.Range("I19").Value2 = 909
$ws.Range("J19").Value2 = 1266
$ws.Range("K19").Value2 = -28.199052132701
$ws.Range("L19").Value2 = -27.454110135674
$ws.Range("M19").Value2 = 31.358381502890
$ws.Range("N19").Value2 = -35.714285714285

# --- Row 20 ---
$ws.Range("C20").Value2 = 32
$ws.Range("D20").Value2 = 48
$ws.Range("E20").Value2 = -33.333333333333
$ws.Range("F20").Value2 = 131
$ws.Range("G20").Value2 = 155
$ws.Range("H20").Value2 = -15.483870967741
$ws.Range("I20").Value2 = 313
$ws.Range("J20").Value2 = 402
$ws.Range("K20").Value2 = -22.139303482587
$ws.Range("L20").Value2 = -22.716049382716
$ws.Range("M20").Value2 = 4.682274247491
$ws.Range("N20").Value2 = -93.055247392944

# --- Row 21 ---
$ws.Range("C21").Value2 = 251
$ws.Range("D21").Value2 = 287
$ws.Range("E21").Value2 = -12.543554006968
$ws.Range("F21").Value2 = 1000
$ws.Range("G21").Value2 = 1237
$ws.Range("H21").Value2 = -19.159256265157
$ws.Range("I21").Value2 = 2369
$ws.Range("J21").Value2 = 2998
$ws.Range("K21").Value2 = -20.980653769179
$ws.Range("L21").Value2 = -19.366916269571
$ws.Range("M21").Value2 = 8.619899128839
$ws.Range("N21").Value2 = -78.945965161749

# --- Row 22 ---
$ws.Range("C22").Value2 = 8
$ws.Range("D22").Value2 = 12
$ws.Range("E22").Value2 = -33.333333333333
$ws.Range("F22").Value2 = 24
$ws.Range("G22").Value2 = 27
$ws.Range("H22").Value2 = -11.111111111111
$ws.Range("I22").Value2 = 51
$ws.Range("J22").Value2 = 54
$ws.Range("K22").Value2 = -5.555555555555
$ws.Range("L22").Value2 = -1.923076923076
$ws.Range("M22").Value2 = 27.5

# --- Row 23 ---
$ws.Range("C23").Value2 = 3
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = 0
$ws.Range("F23").Value2 = 20
$ws.Range("G23").Value2 = 30
$ws.Range("H23").Value2 = -33.333333333333
$ws.Range("I23").Value2 = 44
$ws.Range("J23").Value2 = 46
$ws.Range("K23").Value2 = -4.347826086956
$ws.Range("L23").Value2 = 2.325581395348
$ws.Range("M23").Value2 = 51.724137931034

# --- Row 24 ---
$ws.Range("C24").Value2 = 264
$ws.Range("D24").Value2 = 295
$ws.Range("E24").Value2 = -10.508474576271
$ws.Range("G24").Value2 = 1381
$ws.Range("H24").Value2 = -18.247646632874
$ws.Range("I24").Value2 = 2646
$ws.Range("J24").Value2 = 3252
$ws.Range("K24").Value2 = -18.634686346863
$ws.Range("L24").Value2 = -9.414584046559
$ws.Range("M24").Value2 = 73.508196721311

# --- Row 25 ---
$ws.Range("C25").Value2 = 178
$ws.Range("D25").Value2 = 191
$ws.Range("E25").Value2 = -6.806282722513
$ws.Range("F25").Value2 = 736
$ws.Range("G25").Value2 = 888
$ws.Range("H25").Value2 = -17.117117117117
$ws.Range("I25").Value2 = 1646
$ws.Range("J25").Value2 = 2015
$ws.Range("K25").Value2 = -18.312655086848
$ws.Range("L25").Value2 = 3.327055869428

# --- Row 26 ---
$ws.Range("D26").Value2 = 107
$ws.Range("E26").Value2 = 3.738317757009
$ws.Range("F26").Value2 = 397
$ws.Range("G26").Value2 = 488
$ws.Range("H26").Value2 = -18.647540983606
$ws.Range("I26").Value2 = 928
$ws.Range("J26").Value2 = 1003
$ws.Range("K26").Value2 = -7.477567298105
$ws.Range("L26").Value2 = 6.057142857142
$ws.Range("M26").Value2 = 10.872162485065

# --- Row 27 ---
$ws.Range("C27").Value2 = 4
$ws.Range("D27").Value2 = 5
$ws.Range("E27").Value2 = -20
$ws.Range("F27").Value2 = 14
$ws.Range("G27").Value2 = 20
$ws.Range("H27").Value2 = -30
$ws.Range("I27").Value2 = 46
$ws.Range("J27").Value2 = 45
$ws.Range("K27").Value2 = 2.222222222222
$ws.Range("L27").Value2 = -16.363636363636

# --- Row 28 ---
$ws.Range("C28").Value2 = 10
$ws.Range("D28").Value2 = 13
$ws.Range("E28").Value2 = -23.076923076923
$ws.Range("F28").Value2 = 43
$ws.Range("G28").Value2 = 41
$ws.Range("H28").Value2 = 4.878048780487
$ws.Range("I28").Value2 = 102
$ws.Range("J28").Value2 = 103
$ws.Range("K28").Value2 = -0.970873786407
$ws.Range("L28").Value2 = -8.928571428571

# --- Row 29 ---
# C29: "0" (text) -> 1 (number, style 14)
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value2 = 1
# D29: 2 (number) -> "0" (text, style 13)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4163)
# E29: -100 (number) -> "***.*" (text, style 13)
$ws.Range("N23").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("N23").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("F29").Value2 = 4
$ws.Range("G29").Value2 = 2
$ws.Range("H29").Value2 = 100
$ws.Range("I29").Value2 = 8
$ws.Range("K29").Value2 = 100
$ws.Range("L29").Value2 = -27.272727272727
$ws.Range("M29").Value2 = 14.285714285714
$ws.Range("N29").Value2 = -85.454545454545

# --- Row 30 ---
# C30: "0" (text) -> 1 (number, style 14)
$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value2 = 1
# D30: 1 (number) -> "0" (text, style 13)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4163)
# E30: -100 (number) -> "***.*" (text, style 13)
$ws.Range("N23").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("N23").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("F30").Value2 = 4
$ws.Range("G30").Value2 = 1
$ws.Range("H30").Value2 = 300
$ws.Range("I30").Value2 = 8
$ws.Range("K30").Value2 = 166.666666666667
$ws.Range("L30").Value2 = -20
$ws.Range("M30").Value2 = 33.333333333333
$ws.Range("N30").Value2 = -85.185185185185

# --- Row 31 ---
$ws.Range("C31").Value2 = 1
$ws.Range("D31").Value2 = 2
$ws.Range("F31").Value2 = 5
$ws.Range("H31").Value2 = -28.571428571428
$ws.Range("I31").Value2 = 9
$ws.Range("J31").Value2 = 16
$ws.Range("K31").Value2 = -43.75
$ws.Range("L31").Value2 = -25

# --- Row 33 ---
# C33: 1 (number) -> "0" (text, style 13)
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4163)
# D33: "0" (text) -> 1 (number, style 14)
$ws.Range("I14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value2 = 1
# E33: "***.*" (text) -> -100 (number, style 15)
$ws.Range("K14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value2 = -100
$ws.Range("J33").Value2 = 11
$ws.Range("K33").Value2 = -72.727272727272

$excel.CutCopyMode = $false
